$wb = $excel.ActiveWorkbook

# --- Sheet "FolioPlacement": move selection from F1 to G12 (not the active tab) ---
$wsPlacement = $wb.Worksheets.Item("FolioPlacement")
$wsPlacement.Range("G12").Select()

# --- Sheet "ChildLocationCorrespondences": keep its own selection (J6); it will
#     lose tabSelected once another sheet is activated below ---
# (no explicit change needed here besides losing tabSelected)

# --- Sheet "Removal": restructure columns ---
$ws = $wb.Worksheets.Item("Removal")

# Remove the duplicate ENVIRONMENT_AT_REMOVAL column (old column G); this shifts
# PRIMARY_REASON/ENVIRONMENT_AT_REMOVAL/REMOVED_FROM/REMOVED_BY left by one.
$ws.Columns("G").Delete()

# New trailing header cells (row 1)
$ws.Range("K1").Value = "PERSON_PRESENT_AT_THE_HOME"
$ws.Range("L1").Value = "SECONDARY_REASON_FOR_REMOVAL"
$ws.Range("M1").Value = "PERSON_WHO_LIVED"
$ws.Range("N1").Value = "CHILD_CIRCUMSTANCES"
$ws.Range("O1").Value = "FAMILY_CIRCUMSTANCES"

# Match header formatting (style) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)

# New data cells (row 3)
$ws.Range("I3").Value = "folioPersonName1"
$ws.Range("K3").Value = "n/a"
$ws.Range("L3").Value = "n/a"
$ws.Range("M3").Value = "n/a"
$ws.Range("N3").Value = "n/a"

# Match data-row formatting (style) used by the rest of row 3
$ws.Range("A3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("K3:N3").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Activate "Removal" last so it becomes the selected tab, with L14 selected.
$ws.Activate()
$ws.Range("L14").Select()
